# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders / refreshes the "Periodo Mora" detail table (rows 16-29) on Hoja1
# so that rows are grouped by period (ascending 1904..1910), alternating
# between the two workers (CLAUDIA CECILIA GARCIA MONCRIEFF / PATRICIA
# MARGARITA GARCIA MONCRIEFF), and updates PATRICIA's "Salario Basico"
# (column G) and a couple of "Valor Mora" (column F) amounts to the new,
# corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    @{ Row = 16; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1904"; F = 33125; G = 828116 },
    @{ Row = 17; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1904"; F = 33125; G = 781242 },
    @{ Row = 18; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1905"; F = 33125; G = 828116 },
    @{ Row = 19; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1905"; F = 33125; G = 781242 },
    @{ Row = 20; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1906"; F = 33125; G = 828116 },
    @{ Row = 21; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1906"; F = 31249; G = 781242 },
    @{ Row = 22; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1907"; F = 33125; G = 828116 },
    @{ Row = 23; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1907"; F = 31249; G = 781242 },
    @{ Row = 24; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1908"; F = 33125; G = 828116 },
    @{ Row = 25; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1908"; F = 31249; G = 781242 },
    @{ Row = 26; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1909"; F = 33125; G = 828116 },
    @{ Row = 27; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1909"; F = 31249; G = 781242 },
    @{ Row = 28; B = "CC"; C = "45592819"; D = "CLAUDIA CECILIA GARCIA MONCRIEFF";   E = "1910"; F = 23187; G = 828116 },
    @{ Row = 29; B = "CC"; C = "45592816"; D = "PATRICIA MARGARITA GARCIA MONCRIEFF"; E = "1910"; F = 21874; G = 781242 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = $r.B   # B: Tipo Doc Trabajador
    $ws.Cells.Item($n, 3).Value = $r.C   # C: N Doc Trabajador
    $ws.Cells.Item($n, 4).Value = $r.D   # D: Nombre Trabajador
    $ws.Cells.Item($n, 5).Value = $r.E   # E: Periodo Mora
    $ws.Cells.Item($n, 6).Value = $r.F   # F: Valor Mora
    $ws.Cells.Item($n, 7).Value = $r.G   # G: Salario Basico
}
